{"js": "// Update the \"Total\" row of the evaluation table:\n// the four numeric summary cells change from 6, 4, 1, 1 to 5, 7, 0, 0.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst lastRow = rows.items[rows.items.length - 1];\nconst cells = lastRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\nconst newValues = [\"5\", \"7\", \"0\", \"0\"];\n// The last row has 6 physical cells: Total (merged), 100%, then 4 numeric cells.\nconst numericCells = cells.items.slice(cells.items.length - newValues.length);\nfor (let i = 0; i < numericCells.length; i++) {\n  // Replace just the text content of the cell's range so the existing\n  // paragraph/run formatting (bold, color, rFonts, etc.) is preserved.\n  numericCells[i].getRange().insertText(newValues[i], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Total\" row of the evaluation table:\n# the four numeric summary cells change from 6, 4, 1, 1 to 5, 7, 0, 0.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$lastRow = $t.Rows.Count\n$lastCol = $t.Columns.Count\n\n$newValues = @(\"5\", \"7\", \"0\", \"0\")\n\nfor ($i = 0; $i -lt $newValues.Length; $i++) {\n    $col = $lastCol - $newValues.Length + 1 + $i\n    $t.Cell($lastRow, $col).Range.Text = $newValues[$i]\n}\n"}
